$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three scratch "- d[...]" helper-text cells (M14:M16) are no longer
# needed once the second Levenshtein table (rows 18-28) carries its own
# live formulas, so clear them out. This also drops the now-unused
# shared strings and shrinks the sheet's used range back down to column L.
$ws.Range("M14:M16").ClearContents()

# Fill in the Levenshtein dynamic-programming table for the second word
# pair (rows 18-28, columns D:L), mirroring the pattern already present
# in the first table (rows 1-11).
$ws.Range("D20").Formula = '=MIN(D19+1,C20+1,C19+(IF(D$18=$B20,0,1)))'
$ws.Range("E20:L28").Formula = '=MIN(E19+1,D20+1,D19+(IF(E$18=$B20,0,1)))'
$ws.Range("D21:D28").Formula = '=MIN(D20+1,C21+1,C20+(IF(D$18=$B21,0,1)))'

# Match the author's final selection/scroll state.
[void]$ws.Range("M20").Select()
